$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "DropBag_1" value to Z2:Z6 (this mints a new shared string,
# reused across all five rows, matching the new DropBag_1 shared string).
$ws.Range("Z2").Value = "DropBag_1"
$ws.Range("Z3").Value = "DropBag_1"
$ws.Range("Z4").Value = "DropBag_1"
$ws.Range("Z5").Value = "DropBag_1"
$ws.Range("Z6").Value = "DropBag_1"

# Add a comment on Z1 explaining how to fill the new column.
$comment = $ws.Range("Z1").AddComment("可填入英文分号间隔的奖励包(掉落包和奖励包公用)`n")

# Move the selection / view to match the edited area.
$ws.Range("Z6").Select()
